# Replace the trailing empty paragraph + bookmark paragraph at the end of the
# document with:
#   49. Created tests.
#   50. Added some exceptions of the Position, Cell and PlayerScore classes
#       recommended by some made tests. (Test Driven Development).   [+ _GoBack bookmark]
#   <empty paragraph>
# matching the target diff.

$d = $word.ActiveDocument

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Locate the last two paragraphs of the document: an empty paragraph
# immediately followed by the (also textually empty) paragraph that carries
# the hidden "_GoBack" bookmark. Note: Paragraph.Range.Text includes the
# trailing paragraph-mark character ("`r"), so a genuinely empty paragraph's
# text is just that one character, not the empty string.
$countBefore = $d.Paragraphs.Count
$pEmpty = $d.Paragraphs.Item($countBefore - 1)
$pLast = $d.Paragraphs.Item($countBefore)
if ($pEmpty.Range.Text -ne "`r" -or $pLast.Range.Text -ne "`r") {
    throw "Unexpected document tail: expected the last two paragraphs to be empty (found '$($pEmpty.Range.Text)' / '$($pLast.Range.Text)')."
}

# --- Step 1: split "49. Created tests." off as a brand new paragraph, right
# before the (currently) empty trailing paragraph. Because that paragraph is
# already empty, the old paragraph mark is preserved untouched after the
# split (no stray empty run gets left behind).
$insertion1 = $d.Range($pEmpty.Range.Start, $pEmpty.Range.Start)
$insertion1.InsertBefore("49. Created tests.`r")

# --- Step 2: merge the now-empty paragraph back into the following
# (bookmarked) paragraph by deleting the paragraph mark between them, so the
# bookmark paragraph becomes empty again (still carrying the bookmark).
# The insertion above added one extra paragraph, so the empty one that used
# to be second-to-last is now the new second-to-last index.
$countAfterSplit = $d.Paragraphs.Count
$pEmpty2 = $d.Paragraphs.Item($countAfterSplit - 1)
$mark = $d.Range($pEmpty2.Range.End - 1, $pEmpty2.Range.End)
$mark.Delete()

# --- Step 3: insert the "50. ..." text at the very start of the bookmarked
# paragraph (still a collapsed insertion, so it merges into the same
# paragraph rather than splitting), keeping bookmarkStart/bookmarkEnd in the
# same paragraph as the new text, as in the target.
$countAfterMerge = $d.Paragraphs.Count
$pBookmark = $d.Paragraphs.Item($countAfterMerge)
$insertion2 = $d.Range($pBookmark.Range.Start, $pBookmark.Range.Start)
$insertion2.InsertBefore("50. Added some exceptions of the Position, Cell and PlayerScore classes recommended by some made tests. (Test Driven Development).")

# --- Step 4: append a clean, truly-empty trailing paragraph (no stray run)
# after the bookmarked paragraph by inserting raw paragraph XML at the end
# of the document.
$countAfterText = $d.Paragraphs.Count
$pBookmark2 = $d.Paragraphs.Item($countAfterText)
$endPoint = $d.Range($pBookmark2.Range.End, $pBookmark2.Range.End)
[void]$endPoint.InsertXML("<w:p $wns/>")
